# Atualização de bases das ligas, do dia: 03-05-2024 às 22:15
#
# The underlying match rows got re-paired with their correct ids/results;
# net effect is that, for three row-pairs, everything except the leading
# index column (A) swaps between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(36, 37),
    @(124, 125),
    @(200, 201)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Columns B..AB hold the data to swap; column A is the stable row index.
    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
